$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
}

Replace-Text "2024-05-11 Saturday" "2024-05-12 Sunday"

Replace-Text "592×6=3552" "589×9=5301"
Replace-Text "175×2=350" "354×4=1416"
Replace-Text "313×5=1565" "947×4=3788"
Replace-Text "342×6=2052" "816×4=3264"
Replace-Text "579×4=2316" "719×6=4314"

Replace-Text "558×7=3906" "257×2=514"
Replace-Text "998×2=1996" "391×8=3128"
Replace-Text "667×2=1334" "158×7=1106"
Replace-Text "877×8=7016" "475×8=3800"
Replace-Text "716×6=4296" "781×9=7029"

Replace-Text "536×4=2144" "720×8=5760"
Replace-Text "285×3=855" "362×6=2172"
Replace-Text "997×3=2991" "794×8=6352"
Replace-Text "453×3=1359" "231×6=1386"
Replace-Text "663×2=1326" "286×3=858"

Replace-Text "673×4=2692" "726×8=5808"
Replace-Text "323×9=2907" "575×8=4600"
Replace-Text "902×3=2706" "909×3=2727"
Replace-Text "455×7=3185" "745×4=2980"
Replace-Text "443×8=3544" "564×8=4512"

Replace-Text "692×3=2076" "835×2=1670"
Replace-Text "304×7=2128" "524×5=2620"
Replace-Text "821×9=7389" "452×8=3616"
Replace-Text "732×3=2196" "817×2=1634"
Replace-Text "961×9=8649" "466×8=3728"
